# Automatic update of files.
# Bump the "Förändrad" date column (C) by one day for rows 2-7,
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 45178
